# EPBDS-9540: Support Java Name conversion on Json field name generating
# in SpreadsheetResults. Renames the "Step1"/"Step2"/"SomeStep" JSON-style
# field names referenced in the test/result header cells on Sheet1 to their
# camelCase equivalents: "step1"/"step2"/"someStep".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- mySprFor_text test table (rows 29-30) ---
# Row 29: typed result field names (with ":Integer" suffix)
$ws.Range("D29").Value = '_res_.$Step2["step1"]:Integer'
$ws.Range("E29").Value = '_res_.$Step2["someStep"]:Integer'
$ws.Range("F29").Value = '_res_.$Step2["step2"]:Integer'

# Row 30: plain result field names
$ws.Range("D30").Value = '_res_.$Step2["step1"]'
$ws.Range("E30").Value = '_res_.$Step2["someStep"]'
$ws.Range("F30").Value = '_res_.$Step2["step2"]'

# --- mySpr2d_test test table (rows 45-46) ---
$ws.Range("C45").Value = '_res_.$Step2["step1"]:Integer'
$ws.Range("D45").Value = '_res_.$Step2["step2"]:Integer'
$ws.Range("E45").Value = '_res_.$Step2["step3"]:Integer'

$ws.Range("C46").Value = '_res_.$Step2["step1"]:Integer'
$ws.Range("D46").Value = '_res_.$Step2["step2"]:Integer'
$ws.Range("E46").Value = '_res_.$Step2["step3"]:Integer'

# --- normalize row 48 numeric literals (stored as plain numbers) ---
$ws.Range("B48").Value = 4
$ws.Range("C48").Value = 5
$ws.Range("E48").Value = 7

# --- restore the selection that was active when the workbook was saved ---
$ws.Activate()
$ws.Range("G24").Select()
